# Update countries & provincias Spain
# Applies the daily COVID data refresh to the "Pais" sheet:
#  - updates the "last updated" timestamp label
#  - updates case figures for several countries (rows 4, 17, 95, 96)
#  - Burkina Faso overtakes Uruguay in the ranking, so their two rows swap
#    country names (and therefore their data follows each row position)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp label in row 1
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 03:22"

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 560425
$ws.Cells.Item(4, 3).Value = 23
$ws.Cells.Item(4, 5).Value = 505686

# Row 17 - Brasil
$ws.Cells.Item(17, 2).Value = 22318
$ws.Cells.Item(17, 3).Value = 126
$ws.Cells.Item(17, 5).Value = 20915
$ws.Cells.Item(17, 7).Value = 7
$ws.Cells.Item(17, 8).Value = 1230

# Rows 95/96 - Burkina Faso and Uruguay swap ranking positions
# Row 95 now holds Burkina Faso's figures
$ws.Cells.Item(95, 1).Value = "Burkina Faso"
$ws.Cells.Item(95, 2).Value = 497
$ws.Cells.Item(95, 3).Value = 0
$ws.Cells.Item(95, 4).Value = 161
$ws.Cells.Item(95, 5).Value = 309
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 27

# Row 96 now holds Uruguay's figures
$ws.Cells.Item(96, 1).Value = "Uruguay"
$ws.Cells.Item(96, 2).Value = 480
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 231
$ws.Cells.Item(96, 5).Value = 242
$ws.Cells.Item(96, 6).Value = 16
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 7
